# TC01_Canine_Filter_Breed-Akita.xlsx edit
# Commit: "changed the xlfile akita 01 and created regression suite laxmi_regression1"
#
# The "startup" sheet holds three Neo4j/Cypher queries (one per tab: CasesTab,
# SamplesTab, FilesTab) in column B. The CasesTab query (cell B2) is edited to
# drop its trailing `Cohort` column, since the query no longer returns a
# `co.cohort_description` field:
#
#     RETURN ... coalesce(diag.best_response, '') AS `Response to Treatment`,
#             coalesce(co.cohort_description, '') AS `Cohort`
#   ->
#     RETURN ... coalesce(diag.best_response, '') AS `Response to Treatment`

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Akita']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID`,
        coalesce(s.clinical_study_designation, '') AS `Study Code`,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease`,
        coalesce(demo.patient_age_at_enrollment, '') AS Age,
        coalesce(demo.sex, '') AS Sex,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`

'@

$ws.Range("B2").Value = $casesQuery

# Row heights shrink / shift slightly to match the now-rewrapped cell text.
$ws.Rows.Item(2).RowHeight = 270
$ws.Rows.Item(3).RowHeight = 225
$ws.Rows.Item(4).RowHeight = 255

# The author's final selection lands on the edited cell, B2, instead of C4,
# and the saved view no longer pins a frozen/scrolled top-left cell.
$null = $ws.Range("B2").Select()
